# Applies the commit:
#  1. Slide 16's table switches from the custom "Table_0" style
#     ({F2A7BDFE-ABDE-48BC-B272-8DED8C0F4DB0}) to the built-in style
#     {D8EC0C43-9A5F-4697-9F9B-1D002BA9EAD4}.
#  2. The presentation's (slide master) theme color palette is swapped
#     from the "Integral" palette to the "Office Theme" palette (the
#     font scheme / format scheme are identical between the two themes,
#     only the 12 color-scheme slots actually change).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s = $p.Slides.Item(16)

$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

$tbl = $tableShape.Table
$tbl.ApplyStyle("{D8EC0C43-9A5F-4697-9F9B-1D002BA9EAD4}")

# --- 2. Theme colors --------------------------------------------------
# Order of theme color slots exposed by ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$tcs = $s.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
